# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column (B) with the week's start date
#  - renumber the Week labels from "W01".."W16" to "W1".."W16"
#  - convert the is_holiday_week column to boolean values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (B:I -> C:J).
$ws.Columns("B:B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    # Week label: "W01" -> "W1", ... "W16" stays "W16".
    $weekNum = $i + 1
    $ws.Cells.Item($row, 1).Value = "W$weekNum"

    # Write the date as literal text (leading apostrophe forces Excel to
    # treat it as a string rather than re-interpreting it as a date
    # serial); ClearFormats() drops the quote-prefix style that the
    # apostrophe entry leaves behind so the cell keeps the default style.
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = "'" + $weekStartDates[$i]
    $cell.ClearFormats()

    # is_holiday_week now lives in column J and should be a real boolean.
    $holidayCell = $ws.Cells.Item($row, 10)
    $holidayCell.Value = [bool]($holidayCell.Value2)
}
